$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in strain values for rows 8, 9, 10, 17, 18, 19 (column F)
$ws.Range("F8").Value = "TDY1480"
$ws.Range("F9").Value = "TDY1480"
$ws.Range("F10").Value = "TDY1480"
$ws.Range("F17").Value = "TDY1480"
$ws.Range("F18").Value = "TDY1480"
$ws.Range("F19").Value = "TDY1480"

# Update the selection to match the new active range
$ws.Range("F17:F19").Select()
